$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.427.49"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.232.65"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'578.59"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").Value = "'183.86"
$ws.Range("E6").Value = "  -1.62%  "
$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.229.02"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D11").Value = "'6.59"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "3.786.90"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "'27.55"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "67.449.91"
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D18").Value = "3.239.54"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("D21").Value = "'394.29"
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'71.25"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "'9.57"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").Value = "'5.54"
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").Value = "'22.56"
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("E33").Value = "  -2.99%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("D36").Value = "'161.52"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "'0.801"
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").Value = "'6.47"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").Value = "'2.46"
$ws.Range("E43").Value = "  -6.57%  "
$ws.Range("D44").Value = "'0.0685"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").Value = "2.607.50"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("D47").Value = "'24.67"
$ws.Range("D48").Value = "'334.04"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'6.28"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "  -0.71%  "
